$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "ما اسمك"
$ws.Range("B6").Value = "يونا شات بوت"
$ws.Range("A7").Value = "adsf"
$ws.Range("B7").Value = "asdf"

$ws.Range("B7").Select()
